$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.203.76'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.510.74'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.00'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.506.46'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("E10").Value = '  -0.91%  '

$ws.Range("E11").Value = '  +8.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.583'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.18'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.32%  '

$ws.Range("E14").Value = '  -1.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.074.95'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.28'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '609.75'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.505.19'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.378.71'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.873'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -19.27%  '

$ws.Range("E24").Value = '  -1.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.45'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.09%  '

$ws.Range("E26").Value = '  -4.19%  '

$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("E28").Value = '  -2.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.96'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.98'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.81%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '638.01'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +12.39%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.88'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.71%  '

$ws.Range("E35").Value = '  -4.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.58'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0994'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.40%  '

$ws.Range("E38").Value = '  -0.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0470'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.70'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("E42").Value = '  +1.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0743'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.351.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.92%  '

$ws.Range("E45").Value = '  -5.49%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '32.12'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.84%  '

$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.76'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.28%  '

$ws.Range("E51").Value = '  -0.02%  '
